$d = $word.ActiveDocument

$pairs = @(
    @("33×71=", "60×66="),
    @("57×44=", "25×30="),
    @("27×39=", "39×53="),
    @("37×67=", "90×21="),
    @("50×17=", "15×63="),
    @("18×42=", "33×82="),
    @("22×31=", "81×64="),
    @("64×12=", "52×95="),
    @("14×97=", "68×93="),
    @("75×30=", "46×65="),
    @("75×60=", "24×72="),
    @("70×15=", "94×14="),
    @("29×56=", "21×61="),
    @("60×46=", "21×59="),
    @("62×44=", "24×64="),
    @("22×45=", "66×90="),
    @("88×38=", "39×40="),
    @("36×98=", "25×14="),
    @("62×32=", "54×62="),
    @("41×84=", "97×82="),
    @("72×39=", "32×75="),
    @("28×40=", "77×82="),
    @("45×59=", "59×75="),
    @("67×65=", "72×33="),
    @("87×22=", "62×27=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
